$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out all previously used cells before writing the reshaped table
$ws.UsedRange.ClearContents()

# Header row
$ws.Range("A1").Value = "data"
$ws.Range("B1").Value = "data.name"
$ws.Range("C1").Value = "data.breed"
$ws.Range("D1").Value = "data.age"
$ws.Range("E1").Value = "data.origin"
$ws.Range("F1").Value = "data.origin.country"
$ws.Range("G1").Value = "data.origin.city"

# Row 2 - dog1
$ws.Range("B2").Value = "dog1"
$ws.Range("C2").Value = "dog"
$ws.Range("D2").Value = 2
$ws.Range("F2").Value = "TH"
$ws.Range("G2").Value = "BKK"

# Row 3 - bird1
$ws.Range("B3").Value = "bird1"
$ws.Range("C3").Value = "bird"
$ws.Range("D3").Value = 1
$ws.Range("F3").Value = "TH"
$ws.Range("G3").Value = "BKK"

# Row 4 - cat1
$ws.Range("B4").Value = "cat1"
$ws.Range("C4").Value = "cat"
$ws.Range("D4").Value = 7
$ws.Range("F4").Value = "TH"
$ws.Range("G4").Value = "CNX"

# Row 5 - bird1
$ws.Range("B5").Value = "bird1"
$ws.Range("C5").Value = "bird"
$ws.Range("D5").Value = 1
$ws.Range("F5").Value = "TH"
